# Seventh day - dictionary improvements
# For the listed rows, replace the long Portuguese description in column B
# with the abbreviation already present in column A of the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2,4,5,6,7,8,9,12,16,17,20,21,22,24,26,28,29,30,31,34,35,39,42,43,48,49,50,52,53,54,55,56,59,68,72,73,74,75,76,78,79,80,83,84,92,94,99,107,111,115,119,120,131,132,134,137,139,140,143,144,149,150,151,152,158,162,163,166,171,175,179,180,183,184,189,194,200,203,211,215,216,220,221,225,226,227,230,231,234,235,239,240,242,243,244,245,246,247,248,249,250,251,252,253,265,267,271,280,293)

foreach ($r in $rows) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $aVal
}
